$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group 1 (filters found via repeat finder, short list) - Arial Unicode MS 10pt black
$m1 = $ws.Cells.Item(61,1)
$m1.Style = "Normal"
$m1.Value = "960294-1002017"
$m1.Font.Bold = $false
$m1.Font.Name = "Arial Unicode MS"
$m1.Font.Size = 10
$m1.Font.Color = 0
$m1.Copy()
$ws.Cells.Item(62,1).Value = "1001531-1036138"
$ws.Cells.Item(62,1).PasteSpecial(-4122)
$ws.Cells.Item(63,1).Value = "1098179-1144020"
$ws.Cells.Item(63,1).PasteSpecial(-4122)
$ws.Cells.Item(64,1).Value = "1098179-1144020"
$ws.Cells.Item(64,1).PasteSpecial(-4122)
$ws.Cells.Item(65,1).Value = "1957835-1971770"
$ws.Cells.Item(65,1).PasteSpecial(-4122)
$ws.Cells.Item(66,1).Value = "2728976-2776890"
$ws.Cells.Item(66,1).PasteSpecial(-4122)
$ws.Cells.Item(67,1).Value = "2836642-2885742"
$ws.Cells.Item(67,1).PasteSpecial(-4122)
$ws.Cells.Item(68,1).Value = "2899840-2931190"
$ws.Cells.Item(68,1).PasteSpecial(-4122)
$ws.Cells.Item(69,1).Value = "4417918-4438335"
$ws.Cells.Item(69,1).PasteSpecial(-4122)

# Group 2 (filters found via repeat finder, long list) - Menlo 12pt
$m2 = $ws.Cells.Item(70,1)
$m2.Style = "Normal"
$m2.Value = "815297-815399"
$m2.Font.Bold = $false
$m2.Font.Name = "Menlo"
$m2.Font.Size = 12
$m2.Copy()
$ws.Cells.Item(71,1).Value = "819460-819828"
$ws.Cells.Item(71,1).PasteSpecial(-4122)
$ws.Cells.Item(72,1).Value = "1039690-1039791"
$ws.Cells.Item(72,1).PasteSpecial(-4122)
$ws.Cells.Item(73,1).Value = "1039774-1039869"
$ws.Cells.Item(73,1).PasteSpecial(-4122)
$ws.Cells.Item(74,1).Value = "1224175-1224670"
$ws.Cells.Item(74,1).PasteSpecial(-4122)
$ws.Cells.Item(75,1).Value = "2630737-2631129"
$ws.Cells.Item(75,1).PasteSpecial(-4122)
$ws.Cells.Item(76,1).Value = "2730868-2730950"
$ws.Cells.Item(76,1).PasteSpecial(-4122)
$ws.Cells.Item(77,1).Value = "3077943-3078094"
$ws.Cells.Item(77,1).PasteSpecial(-4122)
$ws.Cells.Item(78,1).Value = "3184549-3184629"
$ws.Cells.Item(78,1).PasteSpecial(-4122)
$ws.Cells.Item(79,1).Value = "3629491-3629571"
$ws.Cells.Item(79,1).PasteSpecial(-4122)
$ws.Cells.Item(80,1).Value = "3682922-3683327"
$ws.Cells.Item(80,1).PasteSpecial(-4122)
$ws.Cells.Item(81,1).Value = "1025020-1025737"
$ws.Cells.Item(81,1).PasteSpecial(-4122)
$ws.Cells.Item(82,1).Value = "1025020-1025738"
$ws.Cells.Item(82,1).PasteSpecial(-4122)
$ws.Cells.Item(83,1).Value = "1224130-1224197"
$ws.Cells.Item(83,1).PasteSpecial(-4122)
$ws.Cells.Item(84,1).Value = "1224682-1224749"
$ws.Cells.Item(84,1).PasteSpecial(-4122)
$ws.Cells.Item(85,1).Value = "2035583-2035659"
$ws.Cells.Item(85,1).PasteSpecial(-4122)
$ws.Cells.Item(86,1).Value = "2045623-2046339"
$ws.Cells.Item(86,1).PasteSpecial(-4122)
$ws.Cells.Item(87,1).Value = "2528936-2529013"
$ws.Cells.Item(87,1).PasteSpecial(-4122)
$ws.Cells.Item(88,1).Value = "2529054-2529131"
$ws.Cells.Item(88,1).PasteSpecial(-4122)
$ws.Cells.Item(89,1).Value = "2531836-2531911"
$ws.Cells.Item(89,1).PasteSpecial(-4122)
$ws.Cells.Item(90,1).Value = "2531957-2532034"
$ws.Cells.Item(90,1).PasteSpecial(-4122)
$ws.Cells.Item(91,1).Value = "2532074-2532149"
$ws.Cells.Item(91,1).PasteSpecial(-4122)
$ws.Cells.Item(92,1).Value = "2532154-2532231"
$ws.Cells.Item(92,1).PasteSpecial(-4122)
$ws.Cells.Item(93,1).Value = "2579219-2579937"
$ws.Cells.Item(93,1).PasteSpecial(-4122)
$ws.Cells.Item(94,1).Value = "2579220-2579937"
$ws.Cells.Item(94,1).PasteSpecial(-4122)
$ws.Cells.Item(95,1).Value = "2579220-2579938"
$ws.Cells.Item(95,1).PasteSpecial(-4122)
$ws.Cells.Item(96,1).Value = "2796440-2796573"
$ws.Cells.Item(96,1).PasteSpecial(-4122)
$ws.Cells.Item(97,1).Value = "294514-294645"
$ws.Cells.Item(97,1).PasteSpecial(-4122)
$ws.Cells.Item(98,1).Value = "294514-294646"
$ws.Cells.Item(98,1).PasteSpecial(-4122)
$ws.Cells.Item(99,1).Value = "294764-294906"
$ws.Cells.Item(99,1).PasteSpecial(-4122)
$ws.Cells.Item(100,1).Value = "294829-294905"
$ws.Cells.Item(100,1).PasteSpecial(-4122)
$ws.Cells.Item(101,1).Value = "2969930-2970008"
$ws.Cells.Item(101,1).PasteSpecial(-4122)
$ws.Cells.Item(102,1).Value = "2970200-2970278"
$ws.Cells.Item(102,1).PasteSpecial(-4122)
$ws.Cells.Item(103,1).Value = "2970478-2970556"
$ws.Cells.Item(103,1).PasteSpecial(-4122)
$ws.Cells.Item(104,1).Value = "304268-304344"
$ws.Cells.Item(104,1).PasteSpecial(-4122)
$ws.Cells.Item(105,1).Value = "3141065-3141147"
$ws.Cells.Item(105,1).PasteSpecial(-4122)
$ws.Cells.Item(106,1).Value = "3141070-3141146"
$ws.Cells.Item(106,1).PasteSpecial(-4122)
$ws.Cells.Item(107,1).Value = "3141171-3141253"
$ws.Cells.Item(107,1).PasteSpecial(-4122)
$ws.Cells.Item(108,1).Value = "3141176-3141252"
$ws.Cells.Item(108,1).PasteSpecial(-4122)
$ws.Cells.Item(109,1).Value = "3194213-3194931"
$ws.Cells.Item(109,1).PasteSpecial(-4122)
$ws.Cells.Item(110,1).Value = "3194214-3194931"
$ws.Cells.Item(110,1).PasteSpecial(-4122)
$ws.Cells.Item(111,1).Value = "3194214-3194932"
$ws.Cells.Item(111,1).PasteSpecial(-4122)
$ws.Cells.Item(112,1).Value = "3316846-3317058"
$ws.Cells.Item(112,1).PasteSpecial(-4122)
$ws.Cells.Item(113,1).Value = "3458510-3458586"
$ws.Cells.Item(113,1).PasteSpecial(-4122)
$ws.Cells.Item(114,1).Value = "3566871-3567002"
$ws.Cells.Item(114,1).PasteSpecial(-4122)
$ws.Cells.Item(115,1).Value = "3566872-3567002"
$ws.Cells.Item(115,1).PasteSpecial(-4122)
$ws.Cells.Item(116,1).Value = "3566872-3567005"
$ws.Cells.Item(116,1).PasteSpecial(-4122)
$ws.Cells.Item(117,1).Value = "3635059-3635775"
$ws.Cells.Item(117,1).PasteSpecial(-4122)
$ws.Cells.Item(118,1).Value = "3682922-3683126"
$ws.Cells.Item(118,1).PasteSpecial(-4122)
$ws.Cells.Item(119,1).Value = "3683122-3683326"
$ws.Cells.Item(119,1).PasteSpecial(-4122)
$ws.Cells.Item(120,1).Value = "4105115-4105248"
$ws.Cells.Item(120,1).PasteSpecial(-4122)
$ws.Cells.Item(121,1).Value = "4105118-4105248"
$ws.Cells.Item(121,1).PasteSpecial(-4122)
$ws.Cells.Item(122,1).Value = "4105118-4105250"
$ws.Cells.Item(122,1).PasteSpecial(-4122)
$ws.Cells.Item(123,1).Value = "4105231-4105373"
$ws.Cells.Item(123,1).PasteSpecial(-4122)
$ws.Cells.Item(124,1).Value = "4105296-4105372"
$ws.Cells.Item(124,1).PasteSpecial(-4122)
$ws.Cells.Item(125,1).Value = "4141148-4141234"
$ws.Cells.Item(125,1).PasteSpecial(-4122)
$ws.Cells.Item(126,1).Value = "4201309-4201439"
$ws.Cells.Item(126,1).PasteSpecial(-4122)
$ws.Cells.Item(127,1).Value = "4356380-4356510"
$ws.Cells.Item(127,1).PasteSpecial(-4122)
$ws.Cells.Item(128,1).Value = "4559283-4560000"
$ws.Cells.Item(128,1).PasteSpecial(-4122)
$ws.Cells.Item(129,1).Value = "4559283-4560001"
$ws.Cells.Item(129,1).PasteSpecial(-4122)
$ws.Cells.Item(130,1).Value = "4596435-4596511"
$ws.Cells.Item(130,1).PasteSpecial(-4122)
$ws.Cells.Item(131,1).Value = "4596435-4596512"
$ws.Cells.Item(131,1).PasteSpecial(-4122)
$ws.Cells.Item(132,1).Value = "4596667-4596743"
$ws.Cells.Item(132,1).PasteSpecial(-4122)
$ws.Cells.Item(133,1).Value = "4596899-4596975"
$ws.Cells.Item(133,1).PasteSpecial(-4122)
$ws.Cells.Item(134,1).Value = "4596899-4596976"
$ws.Cells.Item(134,1).PasteSpecial(-4122)
$ws.Cells.Item(135,1).Value = "4810646-4810734"
$ws.Cells.Item(135,1).PasteSpecial(-4122)
$ws.Cells.Item(136,1).Value = "4810647-4810733"
$ws.Cells.Item(136,1).PasteSpecial(-4122)
$ws.Cells.Item(137,1).Value = "4810878-4810966"
$ws.Cells.Item(137,1).PasteSpecial(-4122)
$ws.Cells.Item(138,1).Value = "4810879-4810965"
$ws.Cells.Item(138,1).PasteSpecial(-4122)
$ws.Cells.Item(139,1).Value = "738249-738327"
$ws.Cells.Item(139,1).PasteSpecial(-4122)
$ws.Cells.Item(140,1).Value = "738341-738416"
$ws.Cells.Item(140,1).PasteSpecial(-4122)
$ws.Cells.Item(141,1).Value = "738451-738526"
$ws.Cells.Item(141,1).PasteSpecial(-4122)
$ws.Cells.Item(142,1).Value = "738642-738720"
$ws.Cells.Item(142,1).PasteSpecial(-4122)
$ws.Cells.Item(143,1).Value = "818781-818858"
$ws.Cells.Item(143,1).PasteSpecial(-4122)
$ws.Cells.Item(144,1).Value = "818988-819063"
$ws.Cells.Item(144,1).PasteSpecial(-4122)
$ws.Cells.Item(145,1).Value = "818988-819065"
$ws.Cells.Item(145,1).PasteSpecial(-4122)
$ws.Cells.Item(146,1).Value = "819067-819144"
$ws.Cells.Item(146,1).PasteSpecial(-4122)
$ws.Cells.Item(147,1).Value = "819192-819269"
$ws.Cells.Item(147,1).PasteSpecial(-4122)
$ws.Cells.Item(148,1).Value = "819402-819479"
$ws.Cells.Item(148,1).PasteSpecial(-4122)
$ws.Cells.Item(149,1).Value = "906556-906624"
$ws.Cells.Item(149,1).PasteSpecial(-4122)
$ws.Cells.Item(150,1).Value = "906727-906939"
$ws.Cells.Item(150,1).PasteSpecial(-4122)
$ws.Cells.Item(151,1).Value = "908736-908804"
$ws.Cells.Item(151,1).PasteSpecial(-4122)

$ws.Range("A61:A69").RowHeight = 17

$ws.Columns.Item(1).ColumnWidth = 19.6666667

$excel.CutCopyMode = 0
Write-Output "done"
